# Apply the commit's changes to the IG-generated StructureDefinition workbook:
#  1. URL metadata value: pythia -> cicada (also mirrored on the Elements sheet's
#     "Fixed Value" for Extension.url, which always echoes the canonical URL).
#  2. Date metadata value updated to the new generation timestamp.
#  3. A new "Jurisdiction" property row (empty value) is inserted into the
#     Metadata table, right after "Contact" and before "Description" -
#     pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- 1 & 2: in-place value edits on the Metadata sheet ---------------------
$meta.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/preferred-interval-status"
$meta.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"

# Mirror the URL change onto the Elements sheet (Extension.url's Fixed Value
# column, R) so both sheets stay consistent with the new canonical URL.
$elements.Cells.Item(5, 18).Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/preferred-interval-status"

# --- 3: insert the new "Jurisdiction" row on the Metadata sheet ------------
$meta.Rows.Item(11).Insert()

# Copy formatting (fill/border/alignment) down from the row above so the new
# row 11 matches the look of the rest of the property table.
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$meta.Cells.Item(11, 1).Value = "Jurisdiction"
$meta.Cells.Item(11, 2).Value = ""
